$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 7

$ws.Range("C8").Value = 29

$ws.Range("C9").Value = 12
$ws.Range("E9").Value = 11

$ws.Range("C10").Value = 3
$ws.Range("F10").Value = 1

$ws.Range("C12").Value = 33

$ws.Range("C13").Value = 60

$ws.Range("C14").Value = 6

$ws.Range("C15").Value = 71

$ws.Range("C22").Value = 26

$ws.Range("C24").Value = 7

$ws.Range("C28").Value = 30

$ws.Range("C29").Value = 7
$ws.Range("F29").Value = 1
